$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.405.48'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.572.38'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.76'
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3767'
$ws.Range("E7").Value = '  +2.78%  '

$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3423'
$ws.Range("E9").Value = '  +1.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.158'
$ws.Range("E10").Value = '  -1.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07649'
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.24'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("E14").Value = '  -1.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.941'
$ws.Range("E15").Value = '  +0.66%  '

$ws.Range("D16").Value = '1.571.91'
$ws.Range("E16").Value = '  +0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001133'
$ws.Range("E17").Value = '  -0.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.09'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06768'
$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.85'
$ws.Range("E21").Value = '  +2.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.214'
$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").Value = '22.402.37'
$ws.Range("E24").Value = '  -0.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.422'
$ws.Range("E25").Value = '  +1.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.704'
$ws.Range("E26").Value = '  -10.63%  '

$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '146.65'
$ws.Range("E28").Value = '  +1.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.034'
$ws.Range("E29").Value = '  +0.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.38'
$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").Value = '1.746.38'
$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.180'
$ws.Range("E32").Value = '  -2.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.011'
$ws.Range("E33").Value = '  +1.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9878'
$ws.Range("E34").Value = '  -6.24%  '

$ws.Range("E35").Value = '  -3.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08602'
$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2313'
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06572'
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.324'
$ws.Range("E40").Value = '  +5.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.475'
$ws.Range("E41").Value = '  -1.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6441'
$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.52'
$ws.Range("E43").Value = '  -3.28%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.05'
$ws.Range("E45").Value = '  -3.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.794'
$ws.Range("E46").Value = '  +0.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6009'
$ws.Range("E47").Value = '  -0.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.301'
$ws.Range("E48").Value = '  +7.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.081'
$ws.Range("E49").Value = '  -2.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.70'
$ws.Range("E50").Value = '  +1.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07328'
$ws.Range("E51").Value = '  +0.41%  '
